$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1756.9375
$ws.Range("I11").Value = 1756.9375
$ws.Range("K11").Value = 1756.9375
$ws.Range("M11").Value = -1616.9375
$ws.Range("H12").Value = 2946
$ws.Range("I12").Value = 2900
$ws.Range("J12").Value = 2992
$ws.Range("K12").Value = 2900
$ws.Range("L12").Value = 2992
$ws.Range("M12").Value = -2730
$ws.Range("N12").Value = -3332
$ws.Range("H28").Value = 1704.2632
$ws.Range("I28").Value = 1243.4445
$ws.Range("K28").Value = 1243.4445
$ws.Range("M28").Value = -758.4445000000001
$ws.Range("H43").Value = 593650.1
$ws.Range("I43").Value = 2850
$ws.Range("K43").Value = 2850
$ws.Range("M43").Value = -2781
$ws.Range("H52").Value = 273.47058
$ws.Range("I52").Value = 175
$ws.Range("K52").Value = 525
$ws.Range("M52").Value = -365
$ws.Range("H86").Value = 149865500
$ws.Range("I86").Value = 285839780
$ws.Range("K86").Value = 285839780
$ws.Range("M86").Value = -285838657
$ws.Range("H89").Value = 149865500
$ws.Range("I89").Value = 285839780
$ws.Range("K89").Value = 1429198900
$ws.Range("M89").Value = -1429193284
$ws.Range("H98").Value = 2529.561
$ws.Range("I98").Value = 2711.7778
$ws.Range("J98").Value = 1217.6
$ws.Range("K98").Value = 2711.7778
$ws.Range("L98").Value = 1217.6
$ws.Range("M98").Value = -1213.7778
$ws.Range("N98").Value = -4213.6
$ws.Range("H122").Value = 2529.561
$ws.Range("I122").Value = 2711.7778
$ws.Range("J122").Value = 1217.6
$ws.Range("K122").Value = 8135.3334
$ws.Range("L122").Value = 3652.8
$ws.Range("M122").Value = -5685.3334
$ws.Range("N122").Value = -8552.799999999999
$ws.Range("H129").Value = 1499.625
$ws.Range("I129").Value = 999.25
$ws.Range("K129").Value = 2997.75
$ws.Range("M129").Value = 2002.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2204312
$ws.Range("I32").Value = 2304682.8
$ws.Range("J32").Value = 21247
$ws.Range("K32").Value = 2304682.8
$ws.Range("L32").Value = 21247
$ws.Range("M32").Value = -2304395.8
$ws.Range("N32").Value = -21821
$ws.Range("H63").Value = 1454.5
$ws.Range("I63").Value = 1454.5
$ws.Range("K63").Value = 1454.5
$ws.Range("M63").Value = -768.5
$ws.Range("H66").Value = 1454.5
$ws.Range("I66").Value = 1454.5
$ws.Range("K66").Value = 7272.5
$ws.Range("M66").Value = -3840.5
$ws.Range("H88").Value = 100000
$ws.Range("J88").Value = 100000
$ws.Range("L88").Value = 100000
$ws.Range("N88").Value = -100812
$ws.Range("H91").Value = 100000
$ws.Range("J91").Value = 100000
$ws.Range("L91").Value = 100000
$ws.Range("N91").Value = -102808
$ws.Range("H122").Value = 2654.4614
$ws.Range("I122").Value = 1963.6111
$ws.Range("J122").Value = 4208.875
$ws.Range("K122").Value = 5890.8333
$ws.Range("L122").Value = 12626.625
$ws.Range("M122").Value = -3440.8333
$ws.Range("N122").Value = -17526.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11118911
$ws.Range("I20").Value = 15156516
$ws.Range("J20").Value = 15497.5
$ws.Range("K20").Value = 15156516
$ws.Range("L20").Value = 15497.5
$ws.Range("M20").Value = -15156269
$ws.Range("N20").Value = -15991.5
$ws.Range("H86").Value = 131752.25
$ws.Range("I86").Value = 149289.86
$ws.Range("J86").Value = 8989
$ws.Range("K86").Value = 149289.86
$ws.Range("L86").Value = 8989
$ws.Range("M86").Value = -148166.86
$ws.Range("N86").Value = -11235
$ws.Range("H89").Value = 131752.25
$ws.Range("I89").Value = 149289.86
$ws.Range("J89").Value = 8989
$ws.Range("K89").Value = 746449.2999999999
$ws.Range("L89").Value = 44945
$ws.Range("M89").Value = -740833.2999999999
$ws.Range("N89").Value = -56177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1142.9231
$ws.Range("I16").Value = 324.7
$ws.Range("K16").Value = 324.7
$ws.Range("M16").Value = -37.69999999999999
$ws.Range("H22").Value = 887.8182
$ws.Range("I22").Value = 876.6
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 876.6
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -526.6
$ws.Range("N22").Value = -1700
$ws.Range("H31").Value = 4542.625
$ws.Range("I31").Value = 2409.7222
$ws.Range("K31").Value = 2409.7222
$ws.Range("M31").Value = -2114.7222
$ws.Range("H34").Value = 4542.625
$ws.Range("I34").Value = 2409.7222
$ws.Range("K34").Value = 2409.7222
$ws.Range("M34").Value = -2207.7222
$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 10000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 10000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -9407
$ws.Range("N42").ClearContents()
$ws.Range("H113").Value = 1142.9231
$ws.Range("I113").Value = 324.7
$ws.Range("K113").Value = 324.7
$ws.Range("M113").Value = 1845.3
$ws.Range("H132").Value = 4754.222
$ws.Range("I132").Value = 3315.2334
$ws.Range("K132").Value = 9945.700199999999
$ws.Range("M132").Value = -7415.700199999999
$ws.Range("H134").Value = 4020.6775
$ws.Range("I134").Value = 1935.3158
$ws.Range("J134").Value = 7322.5
$ws.Range("K134").Value = 5805.9474
$ws.Range("L134").Value = 21967.5
$ws.Range("M134").Value = -3270.9474
$ws.Range("N134").Value = -27037.5
$ws.Range("H141").Value = 261268.73
$ws.Range("J141").Value = 275645.06
$ws.Range("L141").Value = 275645.06
$ws.Range("N141").Value = -286005.06

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 500001950
$ws.Range("I76").Value = 500001950
$ws.Range("K76").Value = 1500005850
$ws.Range("M76").Value = -1500005467
$ws.Range("H79").Value = 500001950
$ws.Range("I79").Value = 500001950
$ws.Range("K79").Value = 1500005850
$ws.Range("M79").Value = -1500004524

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 7338.3335
$ws.Range("I55").Value = 6206
$ws.Range("K55").Value = 6206
$ws.Range("M55").Value = -5879
$ws.Range("H70").Value = 13787.581
$ws.Range("I70").Value = 16282
$ws.Range("J70").Value = 10333.77
$ws.Range("K70").Value = 16282
$ws.Range("L70").Value = 10333.77
$ws.Range("M70").Value = -16012
$ws.Range("N70").Value = -10873.77
$ws.Range("H73").Value = 13787.581
$ws.Range("I73").Value = 16282
$ws.Range("J73").Value = 10333.77
$ws.Range("K73").Value = 16282
$ws.Range("L73").Value = 10333.77
$ws.Range("M73").Value = -15346
$ws.Range("N73").Value = -12205.77
$ws.Range("H80").Value = 94243
$ws.Range("I80").Value = 2735.8
$ws.Range("J80").Value = 170499
$ws.Range("K80").Value = 2735.8
$ws.Range("L80").Value = 170499
$ws.Range("M80").Value = -1737.8
$ws.Range("N80").Value = -172495
$ws.Range("H83").Value = 94243
$ws.Range("I83").Value = 2735.8
$ws.Range("J83").Value = 170499
$ws.Range("K83").Value = 13679
$ws.Range("L83").Value = 852495
$ws.Range("M83").Value = -8687
$ws.Range("N83").Value = -862479
$ws.Range("H97").Value = 1398.0322
$ws.Range("I97").Value = 1158.8462
$ws.Range("J97").Value = 2641.8
$ws.Range("K97").Value = 1158.8462
$ws.Range("L97").Value = 2641.8
$ws.Range("M97").Value = -662.8462
$ws.Range("N97").Value = -3633.8
$ws.Range("H113").Value = 7084.316
$ws.Range("I113").Value = 4531.6665
$ws.Range("K113").Value = 4531.6665
$ws.Range("M113").Value = -2361.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6218.6816
$ws.Range("I7").Value = 5976.091
$ws.Range("K7").Value = 5976.091
$ws.Range("M7").Value = -5864.091
$ws.Range("H40").Value = 7214.75
$ws.Range("I40").Value = 5993
$ws.Range("K40").Value = 5993
$ws.Range("M40").Value = -5857
$ws.Range("H61").Value = 4864.25
$ws.Range("I61").Value = 2053.818
$ws.Range("K61").Value = 2053.818
$ws.Range("M61").Value = -1851.818
$ws.Range("H113").Value = 4864.25
$ws.Range("I113").Value = 2053.818
$ws.Range("K113").Value = 2053.818
$ws.Range("M113").Value = 116.1819999999998
$ws.Range("H122").Value = 3340.1833
$ws.Range("I122").Value = 2864.5386
$ws.Range("J122").Value = 6431.875
$ws.Range("K122").Value = 8593.6158
$ws.Range("L122").Value = 19295.625
$ws.Range("M122").Value = -6143.6158
$ws.Range("N122").Value = -24195.625
$ws.Range("H126").Value = 6218.6816
$ws.Range("I126").Value = 5976.091
$ws.Range("K126").Value = 17928.273
$ws.Range("M126").Value = -15458.273
$ws.Range("H135").Value = 103476.336
$ws.Range("J135").Value = 103476.336
$ws.Range("L135").Value = 103476.336
$ws.Range("N135").Value = -113616.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1339.0605
$ws.Range("I113").Value = 596.48
$ws.Range("K113").Value = 1789.44
$ws.Range("M113").Value = 380.5599999999999
$ws.Range("H126").Value = 1832.375
$ws.Range("I126").Value = 1664.75
$ws.Range("K126").Value = 4994.25
$ws.Range("M126").Value = -2524.25
$ws.Range("H132").Value = 5761.914
$ws.Range("I132").Value = 6863.143
$ws.Range("J132").Value = 4110.0713
$ws.Range("K132").Value = 20589.429
$ws.Range("L132").Value = 12330.2139
$ws.Range("M132").Value = -18059.429
$ws.Range("N132").Value = -17390.2139
$ws.Range("H140").Value = 71499.75
$ws.Range("J140").Value = 71499.75
$ws.Range("L140").Value = 71499.75
$ws.Range("N140").Value = -81859.75
